$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the time-range labels in column C (rows 6-7): shift the 12:xx block
# forward by one slot -- C6 becomes "12:25-12:30" and C7 becomes "12:30-12:35"
# (replacing the old "22:30-22:35" entry, which is no longer needed since the
# following rows already continue the 22:xx sequence).
$ws.Range("C6").Value = "12:25-12:30"
$ws.Range("C7").Value = "12:30-12:35"

# Update the active selection to match the author's final cursor position.
$ws.Range("C13").Select()
